# EPBDS-12566 Add methods Numbers.toString() for number formatting
#
# The Double/Float "Infinity" / "-Infinity" text samples are replaced with
# the shorter "∞" / "-∞" representations now produced by Numbers.toString().
# These values live in cells L12/M12 (miscDouble block) and L24/M24
# (miscFloat block) on Sheet1. They are stored as text (quote-prefixed,
# style index 1) rather than as formula results, so a leading apostrophe is
# used when assigning the new values to keep them as literal text and to
# preserve the existing quote-prefix cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$infinity = [char]0x221E
$negInfinity = "-" + $infinity

$ws.Range("L12").Value = "'" + $infinity
$ws.Range("M12").Value = "'" + $negInfinity

$ws.Range("L24").Value = "'" + $infinity
$ws.Range("M24").Value = "'" + $negInfinity

# Move/restore the active cell selection on the sheet to M23 (was T23).
[void]$ws.Range("M23").Select()
